# fix bug about SalesRecord
#
# A customer ("Daven") bought the hidden/secret "Secret" drink product:
#   - a new Customer record is appended (Customer!A8:G8)
#   - the Drink "Secret" row's stock/sales counters are updated (Drink!D5:E5)
#   - a new SalesRecord row is appended capturing the (buggy, negative
#     userid) purchase (SalesRecord!A7:G7)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Customer sheet: append new customer "Daven"
# ---------------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("Customer")

$customerRow = $wsCustomer.Range("A8:G8")
$customerRow.NumberFormat = "@"
$wsCustomer.Range("A8").Value = "7"
$wsCustomer.Range("B8").Value = "Daven"
$wsCustomer.Range("C8").Value = "SECRET"
$wsCustomer.Range("D8").Value = "0"
$wsCustomer.Range("E8").Value = "0"
$wsCustomer.Range("F8").Value = "123456"
$wsCustomer.Range("G8").Value = "13117826002"
$customerRow.ClearFormats()

# ---------------------------------------------------------------------
# Drink sheet: update stock/sell counters for the "Secret" product (row 5)
# ---------------------------------------------------------------------
$wsDrink = $wb.Worksheets.Item("Drink")

$drinkRow = $wsDrink.Range("D5:E5")
$drinkRow.NumberFormat = "@"
$wsDrink.Range("D5").Value = "996"
$wsDrink.Range("E5").Value = "4"
$drinkRow.ClearFormats()

# ---------------------------------------------------------------------
# SalesRecord sheet: append the new sale record
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("SalesRecord")

$salesRow = $wsSales.Range("A7:G7")
$salesRow.NumberFormat = "@"
$wsSales.Range("A7").Value = "6"
$wsSales.Range("B7").Value = "666"
$wsSales.Range("C7").Value = "Secret"
$wsSales.Range("D7").Value = "1"
$wsSales.Range("E7").Value = "-2"
$wsSales.Range("F7").Value = "0.0"
$wsSales.Range("G7").Value = "2022-12-11"
$salesRow.ClearFormats()
